$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new measurement was recorded for 2026/01/16 (Friday) at hour 3. It sits
# chronologically before the existing 2026/12/29 block, so insert a new row
# at 647 and push everything from the old row 647 onward down by one.
$ws.Rows.Item(647).Insert()

# Column A stores dates as plain text (e.g. "2026/01/16"), not real date
# values. Assigning a date-shaped string straight to .Value makes Excel
# auto-convert it into a date serial, so instead write it as a text
# formula and then collapse the formula down to its literal text result
# via copy / paste-values - this keeps the stored type a plain string
# (matching every other cell in the column) without touching the cell's
# number format/style.
$dateCell = $ws.Cells.Item(647, 1)
$dateCell.Formula = "=""2026/01/16"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(647, 2).Value = "金"
$ws.Cells.Item(647, 3).Value = 3
$ws.Cells.Item(647, 4).Value = 201
